$d = $word.ActiveDocument

# --- Locate the "Client: ..." paragraph and rewrite its text -----------------
$finder = $d.Range(0, $d.Content.End)
$finder.Find.ClearFormatting()
$found = $finder.Find.Execute("Client: Origin Integrated Gas (IG)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $target = $d.Range($finder.Start, $finder.End)

    $newText = "Client: G'Day Gas Distributors Pty Ltd"
    # Direct Range.Text assignment (rather than Find/Replace) keeps the literal
    # straight apostrophe and avoids the smart-quote AutoCorrect/AutoFormat pass.
    $target.Text = $newText

    # Re-acquire the paragraph start so sub-range offsets below are correct.
    $pStart = $target.Start

    # "G'Day" is flagged by Word's proofing engine (a name it doesn't recognize),
    # which is why it ends up isolated in its own run, bracketed by
    # <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>,
    # same as the existing "Wallumbilla" elsewhere in this document. Force that
    # same run isolation here (a transient bookmark guarantees the run break
    # without leaving any residual formatting behind once it is removed).
    $midStart = $pStart + 8   # length of "Client: "
    $midEnd   = $midStart + 5 # length of "G'Day"
    $midRange = $d.Range($midStart, $midEnd)

    $d.Bookmarks.Add("zzTmpSplit", $midRange) | Out-Null
    $d.Bookmarks.Item("zzTmpSplit").Delete()
}

Write-Output $d.Paragraphs.Item(2).Range.Text
